$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C21").Value = 219
$ws.Range("D21").Value = 192
$ws.Range("E21").Value = 27
$ws.Range("F21").Value = 55.01432664756447
